# Apply crypto price/volume updates produced by the scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.645.72'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '2.532.07'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.58%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").Value = '2.530.19'
$ws.Range("E9").Value = '  +0.34%  '
$ws.Range("E10").Value = '  -2.71%  '
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.343'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.62%  '
$ws.Range("D15").Value = '2.983.77'
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000174'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.23%  '
$ws.Range("D17").Value = '67.533.24'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '2.526.45'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '365.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +57.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.15'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.47%  '
$ws.Range("D29").Value = '2.658.97'
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = '0.0₃0942'
$ws.Range("E30").Value = '  -4.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '533.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("E34").Value = '  -4.99%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.128'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.33%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.24'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.50%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("E41").Value = '  -2.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.343'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.996'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '147.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.67'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.40%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.548'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.66%  '
$ws.Range("D50").Value = '0.0₆0273'
$ws.Range("E50").Value = '  -3.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.39%  '
